$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 750
$ws.Range("B3").Value = 400
$ws.Range("B4").Value = 250
$ws.Range("B5").Value = 125
$ws.Range("B6").Value = 28
$ws.Range("B7").Value = 280
$ws.Range("B8").Value = 150
